# Updated cryptos list on Wed Dec 20 07:58:44 UTC 2023 with GitHub Actions
#
# Price (column D) and Volume(1h) (column E) values are stored as TEXT in
# this sheet (coinranking.com scrape). A leading apostrophe forces Excel to
# keep a numeric-looking string (e.g. "255.57") as text instead of
# auto-converting it to a Number; values that already contain a thousands
# separator dot (e.g. "42.979.60") or the percent column (which keeps its
# leading/trailing spaces) are not ambiguous and don't need the prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Bitcoin) ---------------------------------------------------
$ws.Range("D2").Value = "42.979.60"
$ws.Range("E2").Value = "  -0.29%  "

# --- Row 3 (Ethereum) ---------------------------------------------------
$ws.Range("D3").Value = "2.217.41"
$ws.Range("E3").Value = "  -1.45%  "

# --- Row 4 (TetherUSD) ---------------------------------------------------
$ws.Range("E4").Value = "  -0.20%  "

# --- Row 5 (BNB) ---------------------------------------------------
$ws.Range("D5").Value = "'255.57"
$ws.Range("E5").Value = "  +4.23%  "

# --- Row 6 (XRP) ---------------------------------------------------
$ws.Range("D6").Value = "'0.616"
$ws.Range("E6").Value = "  -0.24%  "

# --- Row 7 (Solana) ---------------------------------------------------
$ws.Range("D7").Value = "'76.36"
$ws.Range("E7").Value = "  +0.14%  "

# --- Row 8 (USDC) ---------------------------------------------------
$ws.Range("E8").Value = "  -0.04%  "

# --- Row 9 (Cardano) ---------------------------------------------------
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "  -2.76%  "

# --- Row 10 (Avalanche) ---------------------------------------------------
$ws.Range("D10").Value = "'41.76"
$ws.Range("E10").Value = "  +1.63%  "

# --- Row 11 (Dogecoin) ---------------------------------------------------
$ws.Range("D11").Value = "'0.0911"
$ws.Range("E11").Value = "  -2.66%  "

# --- Row 12 (Polkadot) ---------------------------------------------------
$ws.Range("D12").Value = "'6.94"
$ws.Range("E12").Value = "  -0.91%  "

# --- Row 13 (TRON) ---------------------------------------------------
$ws.Range("D13").Value = "'0.103"
$ws.Range("E13").Value = "  +1.09%  "

# --- Row 14 (WrappedliquidstakedEther2.0) ---------------------------------
$ws.Range("D14").Value = "2.547.75"
$ws.Range("E14").Value = "  -1.71%  "

# --- Row 15 (Chainlink) ---------------------------------------------------
$ws.Range("D15").Value = "'14.47"
$ws.Range("E15").Value = "  -1.25%  "

# --- Row 16 (WrappedEther) ---------------------------------------------------
$ws.Range("D16").Value = "2.206.46"
$ws.Range("E16").Value = "  -2.15%  "

# --- Row 17 (Polygon) ---------------------------------------------------
$ws.Range("D17").Value = "'0.779"
$ws.Range("E17").Value = "  -3.58%  "

# --- Row 18 (WrappedBTC) ---------------------------------------------------
$ws.Range("D18").Value = "42.906.45"
$ws.Range("E18").Value = "  -0.25%  "

# --- Row 19 (ShibaInu) ---------------------------------------------------
$ws.Range("E19").Value = "  -2.87%  "

# --- Row 20 (Litecoin) ---------------------------------------------------
$ws.Range("D20").Value = "'71.36"
$ws.Range("E20").Value = "  +0.27%  "

# --- Row 21 (Uniswap) ---------------------------------------------------
$ws.Range("D21").Value = "'5.96"
$ws.Range("E21").Value = "  -0.68%  "

# --- Row 22 / Row 23 swap: ImmutableX <-> BitcoinCash -----------------------

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'230.34"
$ws.Range("E22").Value = "  -0.27%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'2.20"
$ws.Range("E23").Value = "  -0.70%  "

# --- Row 24 (InternetComputer) ---------------------------------------------------
$ws.Range("D24").Value = "'9.26"
$ws.Range("E24").Value = "  -8.49%  "

# --- Row 25 (Dai) ---------------------------------------------------
$ws.Range("E25").Value = "  -0.10%  "

# --- Row 26 / Row 27 swap: Cosmos <-> InjectiveProtocol ---------------------

$ws.Range("B26").Value = "InjectiveProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D26").Value = "'41.52"
$ws.Range("E26").Value = "  +4.84%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'10.72"
$ws.Range("E27").Value = "  -1.95%  "

# --- Row 28 (WEMIXToken) ---------------------------------------------------
$ws.Range("E28").Value = "  -3.39%  "

# --- Row 29 (Toncoin) ---------------------------------------------------
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  +0.20%  "

# --- Row 30 (PancakeSwap) ---------------------------------------------------
$ws.Range("E30").Value = "  -2.96%  "

# --- Row 31 (Monero) ---------------------------------------------------
$ws.Range("D31").Value = "'174.01"
$ws.Range("E31").Value = "  +0.16%  "

# --- Row 32 (EthereumClassic) ---------------------------------------------------
$ws.Range("D32").Value = "'20.30"
$ws.Range("E32").Value = "  -0.35%  "

# --- Row 33 (Hedera) ---------------------------------------------------
$ws.Range("D33").Value = "'0.0862"
$ws.Range("E33").Value = "  +8.05%  "

# --- Row 34 (Filecoin) ---------------------------------------------------
$ws.Range("D34").Value = "'5.23"
$ws.Range("E34").Value = "  -2.44%  "

# --- Row 35 (Stellar) ---------------------------------------------------
$ws.Range("E35").Value = "  -0.94%  "

# --- Row 36 (VeChain) ---------------------------------------------------
$ws.Range("D36").Value = "'0.0354"
$ws.Range("E36").Value = "  +5.71%  "

# --- Row 37 (Kaspa) ---------------------------------------------------
$ws.Range("E37").Value = "  -3.57%  "

# --- Row 38 (RenderToken) ---------------------------------------------------
$ws.Range("E38").Value = "  -1.45%  "

# --- Row 39 (Celestia) ---------------------------------------------------
$ws.Range("D39").Value = "'12.58"
$ws.Range("E39").Value = "  -3.50%  "

# --- Row 40 (LidoDAOToken) ---------------------------------------------------
$ws.Range("E40").Value = "  -2.01%  "

# --- Row 41 (NEARProtocol) ---------------------------------------------------
$ws.Range("D41").Value = "'2.75"
$ws.Range("E41").Value = "  +14.65%  "

# --- Row 42 (THORChain) ---------------------------------------------------
$ws.Range("E42").Value = "  -4.50%  "

# --- Row 43 / Row 44 swap: Algorand <-> MultiversX --------------------------

$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").Value = "'60.33"
$ws.Range("E43").Value = "  +0.34%  "

$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").Value = "'0.198"
$ws.Range("E44").Value = "  -2.82%  "

# --- Row 45 (Aave) ---------------------------------------------------
$ws.Range("D45").Value = "'102.22"
$ws.Range("E45").Value = "  -3.60%  "

# --- Row 46 / Row 47 swap: Cronos <-> FraxShare ------------------------------

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'8.31"
$ws.Range("E46").Value = "  -4.52%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0977"
$ws.Range("E47").Value = "  -2.22%  "

# --- Row 48 (WOONetwork) ---------------------------------------------------
$ws.Range("D48").Value = "'0.451"
$ws.Range("E48").Value = "  -2.65%  "

# --- Row 49 (ARBITRUM) ---------------------------------------------------
$ws.Range("E49").Value = "  +0.87%  "

# --- Row 50 (TrustWalletToken) ---------------------------------------------------
$ws.Range("E50").Value = "  -1.51%  "

# --- Row 51 (RocketPoolETH) ---------------------------------------------------
$ws.Range("D51").Value = "2.438.55"
$ws.Range("E51").Value = "  -0.97%  "
